$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-18.91755273140541,-0.3347014285142614,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,0.5680663104253546,-18.91755273140541,-18.91755273140541),
    @(-18.91755273140541,-0.2687703528589985,0.3293757057756997,-18.91755273140541,3.929284507907307,-18.91755273140541,1.729400552336519,-18.91755273140541,2.672981781426024,-18.91755273140541),
    @(-18.91755273140541,0.4657486658369932,-18.91755273140541,-18.91755273140541,-18.91755273140541,3.509644747600384,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,1.928091216361349,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541),
    @(4.321925327004261,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,0.5242911747346429,-18.91755273140541,2.010535469558099),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,1.631145965166104,-18.91755273140541,1.556621578633006,-18.91755273140541,-18.91755273140541,-18.91755273140541,1.447809629662431),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,2.108833598644674,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,1.037712626399033,2.529849880017292),
    @(-18.91755273140541,-18.91755273140541,1.583702322105744,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,1.789631965161362),
    @(-18.91755273140541,-18.91755273140541,-0.01018965071264297,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,2.209855380666116,-18.91755273140541),
    @(-18.91755273140541,0.542889166010005,-0.07819821927835682,-18.91755273140541,-18.91755273140541,-18.91755273140541,2.177656900503785,0.3641079419843613,1.834289581007948,-18.91755273140541),
    @(-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,-18.91755273140541,2.118245801457239,-0.2652644178773622,1.755867635492492,-18.91755273140541),
    @(-18.91755273140541,-18.91755273140541,2.958664796692303,-18.91755273140541,-18.91755273140541,-18.91755273140541,1.826342367191892,1.612096097127433,-18.91755273140541,-18.91755273140541),
    @(-18.91755273140541,2.870619546275878,2.592607296592946,-18.91755273140541,2.252570991857621,-18.91755273140541,1.070822689289502,3.573268038615133,-18.91755273140541,2.003722363284824),
    @(-18.91755273140541,3.040065500721203,-18.91755273140541,3.134786125729074,-18.91755273140541,2.503024857060449,1.118514875664911,-18.91755273140541,-18.91755273140541,-18.91755273140541)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
